# Append rows 111-116 to the monitor_price log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2024-09-29 03:54:48", "monitor_price", "https://example.com/product", "$199.99", "2024-09-29", "03:54:48"),
    @("2024-09-29 03:54:48", "monitor_price", "invalid_url", "Error fetching price: Invalid URL", "2024-09-29", "03:54:48"),
    @("2024-09-29 03:54:50", "monitor_price", "https://example.com/product", "100 USD", "2024-09-29", "03:54:50"),
    @("2024-09-29 03:56:05", "monitor_price", "https://example.com/product", "$199.99", "2024-09-29", "03:56:05"),
    @("2024-09-29 03:56:05", "monitor_price", "invalid_url", "Error fetching price: Invalid URL", "2024-09-29", "03:56:05"),
    @("2024-09-29 03:56:07", "monitor_price", "https://example.com/product", "100 USD", "2024-09-29", "03:56:07")
)

$startRow = 111

# Force text interpretation so values like "$199.99" / "2024-09-29" are not
# auto-coerced into Number/Date by the Value setter, matching the source
# workbook where every cell is a literal inline string.
$endRow = $startRow + $rows.Count - 1
$fillRange = $ws.Range("A$startRow`:F$endRow")
$fillRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowData[$col - 1]
    }
}

# Revert the temporary text format so the new cells keep the workbook's
# default (unstyled) look, same as the existing data rows.
$fillRange.ClearFormats()
